# Final cleanup pass on Starttab.xlsx:
#  - the icons sheet's logo filename lost its stray "3" (BCN_Logo3.png -> BCN_Logo.png)
#  - the live selection moved on to B10:B11 on the icons sheet

$wb = $excel.ActiveWorkbook

$icons = $wb.Worksheets.Item("icons")
$icons.Activate()

# Every cell that held the old "BCN_Logo3.png" shared string becomes "BCN_Logo.png".
$icons.Range("B2").Value = "BCN_Logo.png"
$icons.Range("B3").Value = "BCN_Logo.png"

# Selection on the active sheet moves to B10:B11.
$icons.Range("B10:B11").Select()
